$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct test data: "SNP07" -> "S-NP07"
$ws.Range("D11").Value = "S-NP07"
$ws.Range("D14").Value = "S-NP07"

# Move the active selection to D4
$ws.Range("D4").Select()
